# Insert a new data row at row 6, pushing the existing rows 6..124 down to
# 7..125 (this also grows the sheet's used range from A1:T124 to A1:T125,
# and the former last row, old row 124, ends up at row 125 - exactly what
# the target workbook needs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new record.
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44921
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101001
$ws.Range("J6").Value = "Arándano (blue)"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 2000
$ws.Range("O6").Value = 2000
$ws.Range("P6").Value = 2000
$ws.Range("Q6").Value = "$/kilo"
$ws.Range("R6").Value = "Región del Maule"
$ws.Range("S6").Value = 2000
$ws.Range("T6").Value = 1
